$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 22; this shifts the existing
# rows 22-29 down to 23-30 (carrying all of their data/formatting along).
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new weekly record.
$ws.Range("A22").Value = 7
$ws.Range("B22").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C22").Value = "Ñuble"
$ws.Range("D22").Value = 44523
$ws.Range("E22").Value = 16
$ws.Range("F22").Value = 100112026
$ws.Range("G22").Value = "Haba"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 7000
$ws.Range("M22").Value = 6500
$ws.Range("N22").Value = "$/saco 25 kilos"
$ws.Range("O22").Value = "Provincia de Diguillín"
$ws.Range("P22").Value = 260
$ws.Range("Q22").Value = 25
$ws.Range("R22").Value = "Hortaliza"
